$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New block of rows for subject S6 / robert, session 2013-06-28, runs 2-9
# (mirrors the existing per-session blocks already in the sheet)

$rows = @(
    @("S6", "robert", "28/06/2013", "2013-06-28-robert", "2013-06-28-10-53-09-run2", 1),
    @("S6", "robert", "28/06/2013", "2013-06-28-robert", "2013-06-28-11-01-33-run3", 2),
    @("S6", "robert", "28/06/2013", "2013-06-28-robert", "2013-06-28-11-08-20-run4", 3),
    @("S6", "robert", "28/06/2013", "2013-06-28-robert", "2013-06-28-11-19-58-run5", 4),
    @("S6", "robert", "28/06/2013", "2013-06-28-robert", "2013-06-28-11-27-28-run6", 5),
    @("S6", "robert", "28/06/2013", "2013-06-28-robert", "2013-06-28-11-36-41-run7", 6),
    @("S6", "robert", "28/06/2013", "2013-06-28-robert", "2013-06-28-11-43-24-run8", 7),
    @("S6", "robert", "28/06/2013", "2013-06-28-robert", "2013-06-28-11-50-00-run9", 8)
)

$startRow = 42
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 1).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 2).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 3).Value = $data[2]

    $ws.Cells.Item($r, 4).Value = $data[3]

    $ws.Cells.Item($r, 5).Value = $data[4]

    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 6).HorizontalAlignment = -4108
}

# match the cursor/selection state left behind in the saved file
$ws.Range("F50").Select()

Write-Output "done"
